# Updated cryptos list (price + 1h volume change refresh) — GitHub Actions data pull
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = "25.824.67"
$ws.Cells.Item(2, 5).Value = "  -0.05%  "

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = "1.635.88"
$ws.Cells.Item(3, 5).Value = "  +0.02%  "

# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = "  +0.04%  "

# Row 5: BNB
$ws.Cells.Item(5, 4).Value = "'215.33"
$ws.Cells.Item(5, 5).Value = "  -0.23%  "

# Row 6: XRP
$ws.Cells.Item(6, 5).Value = "  -0.60%  "

# Row 7: USDC
$ws.Cells.Item(7, 5).Value = "  +0.03%  "

# Row 8: Cardano
$ws.Cells.Item(8, 5).Value = "  -0.20%  "

# Row 9: Dogecoin
$ws.Cells.Item(9, 4).Value = "'0.0643"
$ws.Cells.Item(9, 5).Value = "  -0.15%  "

# Row 10: Solana
$ws.Cells.Item(10, 4).Value = "'19.89"
$ws.Cells.Item(10, 5).Value = "  +1.80%  "

# Row 11: TRON
$ws.Cells.Item(11, 4).Value = "'0.0783"
$ws.Cells.Item(11, 5).Value = "  +0.47%  "

# Row 12: Polkadot
$ws.Cells.Item(12, 5).Value = "  -0.78%  "

# Row 13: WrappedEther
$ws.Cells.Item(13, 4).Value = "1.643.78"
$ws.Cells.Item(13, 5).Value = "  +0.61%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Cells.Item(14, 4).Value = "1.860.65"
$ws.Cells.Item(14, 5).Value = "  +0.13%  "

# Row 15: Polygon
$ws.Cells.Item(15, 5).Value = "  -0.61%  "

# Row 16: ShibaInu
$ws.Cells.Item(16, 5).Value = "  +1.44%  "

# Row 17: Litecoin
$ws.Cells.Item(17, 4).Value = "'63.12"
$ws.Cells.Item(17, 5).Value = "  -0.04%  "

# Row 18: WrappedBTC
$ws.Cells.Item(18, 4).Value = "25.828.30"
$ws.Cells.Item(18, 5).Value = "  -0.02%  "

# Row 19: Dai
$ws.Cells.Item(19, 5).Value = "  -0.04%  "

# Row 20: Uniswap
$ws.Cells.Item(20, 5).Value = "  +1.87%  "

# Row 21: BitcoinCash
$ws.Cells.Item(21, 4).Value = "'194.28"
$ws.Cells.Item(21, 5).Value = "  -0.17%  "

# Row 22: Avalanche
$ws.Cells.Item(22, 4).Value = "'9.95"
$ws.Cells.Item(22, 5).Value = "  +0.94%  "

# Row 23: Chainlink
$ws.Cells.Item(23, 5).Value = "  +1.52%  "

# Row 24: BinanceUSD
$ws.Cells.Item(24, 5).Value = "  +0.03%  "

# Row 25: Toncoin
$ws.Cells.Item(25, 4).Value = "'1.76"
$ws.Cells.Item(25, 5).Value = "  -1.55%  "

# Row 26: Monero
$ws.Cells.Item(26, 4).Value = "'139.45"
$ws.Cells.Item(26, 5).Value = "  -0.80%  "

# Row 27: Stellar
$ws.Cells.Item(27, 5).Value = "  -5.27%  "

# Row 28: Cosmos
$ws.Cells.Item(28, 4).Value = "'6.84"
$ws.Cells.Item(28, 5).Value = "  +1.21%  "

# Row 29: EthereumClassic
$ws.Cells.Item(29, 4).Value = "'15.55"
$ws.Cells.Item(29, 5).Value = "  +0.82%  "

# Row 30: PancakeSwap
$ws.Cells.Item(30, 5).Value = "  +0.28%  "

# Row 31: Hedera
$ws.Cells.Item(31, 4).Value = "'0.0496"
$ws.Cells.Item(31, 5).Value = "  +1.50%  "

# Row 32: InternetComputer(DFINITY)
$ws.Cells.Item(32, 5).Value = "  +1.17%  "

# Row 33: Filecoin
$ws.Cells.Item(33, 5).Value = "  +1.28%  "

# Row 34: LidoDAOToken
$ws.Cells.Item(34, 5).Value = "  +2.31%  "

# Row 35: HuobiToken
$ws.Cells.Item(35, 5).Value = "  +0.42%  "

# Row 36: ARBITRUM
$ws.Cells.Item(36, 5).Value = "  +0.47%  "

# Row 37: MXToken
$ws.Cells.Item(37, 5).Value = "  +0.17%  "

# Row 38: ImmutableX
$ws.Cells.Item(38, 5).Value = "  +0.56%  "

# Row 39: Maker
$ws.Cells.Item(39, 4).Value = "1.111.87"
$ws.Cells.Item(39, 5).Value = "  -1.59%  "

# Row 40: VeChain
$ws.Cells.Item(40, 4).Value = "'0.0156"
$ws.Cells.Item(40, 5).Value = "  +0.30%  "

# Row 41: PaxDollar
$ws.Cells.Item(41, 5).Value = "  +0.65%  "

# Row 42: FraxShare
$ws.Cells.Item(42, 5).Value = "  +0.81%  "

# Row 43: Quant/TrustWalletToken (rows swapped)
$ws.Cells.Item(43, 2).Value = "TrustWalletToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(43, 4).Value = "'0.804"
$ws.Cells.Item(43, 5).Value = "  +0.51%  "

# Row 44: TrustWalletToken/Quant (rows swapped)
$ws.Cells.Item(44, 2).Value = "Quant"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(44, 4).Value = "'99.42"
$ws.Cells.Item(44, 5).Value = "  +2.19%  "

# Row 45: BabyDogeCoin
$ws.Cells.Item(45, 4).Value = "0.0₆0111"
$ws.Cells.Item(45, 5).Value = "  -0.40%  "

# Row 46: SynthetixNetwork
$ws.Cells.Item(46, 4).Value = "'2.56"
$ws.Cells.Item(46, 5).Value = "  +14.18%  "

# Row 47: Aave
$ws.Cells.Item(47, 4).Value = "'55.53"
$ws.Cells.Item(47, 5).Value = "  +0.03%  "

# Row 48: Mantle
$ws.Cells.Item(48, 5).Value = "  -5.62%  "

# Row 49: EnergySwap
$ws.Cells.Item(49, 4).Value = "'7.69"
$ws.Cells.Item(49, 5).Value = "  +0.14%  "

# Row 50: Cronos
$ws.Cells.Item(50, 5).Value = "  -0.45%  "

# Row 51: Frax
$ws.Cells.Item(51, 5).Value = "  +0.43%  "
